$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QRVO")

# Row 6: Change in inventories
$ws.Range("B6").Value = 4240000.0
$ws.Range("C6").Value = 10222000.0
$ws.Range("D6").Value = -3562000.0
$ws.Range("E6").Value = 10252000.0
$ws.Range("F6").Value = -2201000.0

# Row 7: Change in payables and accrued liability
$ws.Range("B7").Value = 282597000.0
$ws.Range("C7").Value = 309154000.0
$ws.Range("D7").Value = 245830000.0
$ws.Range("E7").Value = 288723000.0
$ws.Range("F7").Value = 62801000.0
